$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Note sur 4" header in I1, centered (new style picks up horizontal="center")
$ws.Range("I1").Value = "Note sur 4"
$ws.Range("I1").HorizontalAlignment = -4108   # xlCenter

# I5 gets the same centered alignment applied to its existing value
$ws.Range("I5").HorizontalAlignment = -4108   # xlCenter

# Mirror column A (names) into new column J for rows 2-19
$ws.Range("A2:A19").Copy($ws.Range("J2"))

# Column I is narrower now that it only holds short numeric grades
$ws.Columns("I").ColumnWidth = 9.43

# Match the selection left behind by the edit
[void]$ws.Range("J8").Select()
